# Auto-generated script applying the 'Actualizacion 10 de Mayo' diff
# Updates statistical values (Aprobados, Reprobados, Por_Apro, Por_Repro, Promedio, Blancos, Por_Blan)
# across the three sheets: 1er Parcial, 2o Parcial, 3er Parcial

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Cells.Item(2, 5).Value = 33
$ws.Cells.Item(2, 6).Value = 6
$ws.Cells.Item(2, 7).Value = 84.62
$ws.Cells.Item(2, 8).Value = 15.38
$ws.Cells.Item(2, 9).Value = 9.300000000000001
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(2, 11).Value = 15.38
$ws.Cells.Item(9, 5).Value = 29
$ws.Cells.Item(9, 6).Value = 10
$ws.Cells.Item(9, 7).Value = 74.36
$ws.Cells.Item(9, 8).Value = 25.64
$ws.Cells.Item(9, 9).Value = 7.7
$ws.Cells.Item(9, 10).Value = 10
$ws.Cells.Item(9, 11).Value = 25.64
$ws.Cells.Item(18, 5).Value = 29
$ws.Cells.Item(18, 6).Value = 8
$ws.Cells.Item(18, 7).Value = 78.38
$ws.Cells.Item(18, 8).Value = 21.62
$ws.Cells.Item(18, 9).Value = 7.3
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 20
$ws.Cells.Item(19, 7).Value = 44.44
$ws.Cells.Item(19, 8).Value = 55.56
$ws.Cells.Item(19, 9).Value = 6.3
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(20, 5).Value = 13
$ws.Cells.Item(20, 6).Value = 16
$ws.Cells.Item(20, 7).Value = 44.83
$ws.Cells.Item(20, 8).Value = 55.17
$ws.Cells.Item(20, 9).Value = 6.1
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(21, 9).Value = 7.8
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(22, 5).Value = 16
$ws.Cells.Item(22, 6).Value = 6
$ws.Cells.Item(22, 7).Value = 72.73
$ws.Cells.Item(22, 8).Value = 27.27
$ws.Cells.Item(22, 9).Value = 7.5
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Cells.Item(2, 5).Value = 26
$ws.Cells.Item(2, 6).Value = 13
$ws.Cells.Item(2, 7).Value = 66.67
$ws.Cells.Item(2, 8).Value = 33.33
$ws.Cells.Item(2, 9).Value = 9.800000000000001
$ws.Cells.Item(2, 10).Value = 13
$ws.Cells.Item(2, 11).Value = 33.33
$ws.Cells.Item(3, 5).Value = 27
$ws.Cells.Item(3, 6).Value = 8
$ws.Cells.Item(3, 7).Value = 77.14
$ws.Cells.Item(3, 8).Value = 22.86
$ws.Cells.Item(3, 10).Value = 8
$ws.Cells.Item(3, 11).Value = 22.86
$ws.Cells.Item(4, 5).Value = 27
$ws.Cells.Item(4, 6).Value = 8
$ws.Cells.Item(4, 7).Value = 77.14
$ws.Cells.Item(4, 8).Value = 22.86
$ws.Cells.Item(4, 9).Value = 9.6
$ws.Cells.Item(4, 10).Value = 8
$ws.Cells.Item(4, 11).Value = 22.86
$ws.Cells.Item(9, 5).Value = 22
$ws.Cells.Item(9, 6).Value = 17
$ws.Cells.Item(9, 7).Value = 56.41
$ws.Cells.Item(9, 8).Value = 43.59
$ws.Cells.Item(9, 9).Value = 6.9
$ws.Cells.Item(9, 10).Value = 17
$ws.Cells.Item(9, 11).Value = 43.59
$ws.Cells.Item(15, 5).Value = 39
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 100
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 8.199999999999999
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(16, 5).Value = 39
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 100
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 8.199999999999999
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(17, 5).Value = 37
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 100
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 6.9
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(18, 5).Value = 18
$ws.Cells.Item(18, 6).Value = 19
$ws.Cells.Item(18, 7).Value = 48.65
$ws.Cells.Item(18, 8).Value = 51.35
$ws.Cells.Item(18, 9).Value = 8.300000000000001
$ws.Cells.Item(18, 10).Value = 19
$ws.Cells.Item(18, 11).Value = 51.35
$ws.Cells.Item(19, 5).Value = 7
$ws.Cells.Item(19, 6).Value = 29
$ws.Cells.Item(19, 7).Value = 19.44
$ws.Cells.Item(19, 8).Value = 80.56
$ws.Cells.Item(19, 9).Value = 9
$ws.Cells.Item(19, 10).Value = 29
$ws.Cells.Item(19, 11).Value = 80.56
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 21
$ws.Cells.Item(20, 7).Value = 27.59
$ws.Cells.Item(20, 8).Value = 72.41
$ws.Cells.Item(20, 9).Value = 8.5
$ws.Cells.Item(20, 10).Value = 21
$ws.Cells.Item(20, 11).Value = 72.41
$ws.Cells.Item(21, 5).Value = 27
$ws.Cells.Item(21, 6).Value = 9
$ws.Cells.Item(21, 7).Value = 75
$ws.Cells.Item(21, 8).Value = 25
$ws.Cells.Item(21, 9).Value = 8.9
$ws.Cells.Item(21, 10).Value = 9
$ws.Cells.Item(21, 11).Value = 25
$ws.Cells.Item(22, 5).Value = 10
$ws.Cells.Item(22, 6).Value = 12
$ws.Cells.Item(22, 7).Value = 45.45
$ws.Cells.Item(22, 8).Value = 54.55
$ws.Cells.Item(22, 9).Value = 8.9
$ws.Cells.Item(22, 10).Value = 12
$ws.Cells.Item(22, 11).Value = 54.55
$ws = $wb.Worksheets.Item("3er Parcial")
$ws.Cells.Item(2, 5).Value = 33
$ws.Cells.Item(2, 6).Value = 6
$ws.Cells.Item(2, 7).Value = 84.62
$ws.Cells.Item(2, 8).Value = 15.38
$ws.Cells.Item(2, 9).Value = 9.4
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(2, 11).Value = 15.38
$ws.Cells.Item(4, 9).Value = 9.5
$ws.Cells.Item(9, 5).Value = 29
$ws.Cells.Item(9, 6).Value = 10
$ws.Cells.Item(9, 7).Value = 74.36
$ws.Cells.Item(9, 8).Value = 25.64
$ws.Cells.Item(9, 9).Value = 7.7
$ws.Cells.Item(9, 10).Value = 10
$ws.Cells.Item(9, 11).Value = 25.64
$ws.Cells.Item(15, 9).Value = 8.9
$ws.Cells.Item(16, 9).Value = 8.9
$ws.Cells.Item(17, 9).Value = 7.9
$ws.Cells.Item(18, 5).Value = 29
$ws.Cells.Item(18, 6).Value = 8
$ws.Cells.Item(18, 7).Value = 78.38
$ws.Cells.Item(18, 8).Value = 21.62
$ws.Cells.Item(18, 9).Value = 7.4
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 20
$ws.Cells.Item(19, 7).Value = 44.44
$ws.Cells.Item(19, 8).Value = 55.56
$ws.Cells.Item(19, 9).Value = 6.3
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(20, 5).Value = 14
$ws.Cells.Item(20, 6).Value = 15
$ws.Cells.Item(20, 7).Value = 48.28
$ws.Cells.Item(20, 8).Value = 51.72
$ws.Cells.Item(20, 9).Value = 6.3
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(21, 5).Value = 30
$ws.Cells.Item(21, 6).Value = 6
$ws.Cells.Item(21, 7).Value = 83.33
$ws.Cells.Item(21, 8).Value = 16.67
$ws.Cells.Item(21, 9).Value = 8.199999999999999
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(22, 5).Value = 16
$ws.Cells.Item(22, 6).Value = 6
$ws.Cells.Item(22, 7).Value = 72.73
$ws.Cells.Item(22, 8).Value = 27.27
$ws.Cells.Item(22, 9).Value = 7.6
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0

Write-Host "Applied Electricidad - Estadisticos 2020 update (181 cells across 3 sheets)"
